# Update need_to_buy.xlsx values (refresh from R) on Sheet 1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$ws.Range("B2").Value = 14174.6402587165
$ws.Range("C2").Value = 14030.7171968703
$ws.Range("D2").Value = 23859.86
$ws.Range("E2").Value = 9343.54677284558
$ws.Range("F2").Value = -20.2331679285048

$ws.Range("B3").Value = 13458.2767378493
$ws.Range("C3").Value = 12092.7609666455
$ws.Range("E3").Value = 9183.80063892639
$ws.Range("F3").Value = 383.362566898828

$ws.Range("B4").Value = 5454.13158704677
$ws.Range("C4").Value = 8680.33870027015
$ws.Range("E4").Value = 8628.68022479173
$ws.Range("F4").Value = 218.048288544245

$ws.Range("B5").Value = 5395.09658713259
$ws.Range("C5").Value = 8965.62440463466
$ws.Range("E5").Value = 8706.90910918547
$ws.Range("F5").Value = 233.194729742506

$ws.Range("B6").Value = 13159.80815732
$ws.Range("C6").Value = 13224.9542198926
$ws.Range("E6").Value = 8920.98876673794
$ws.Range("F6").Value = 419.586791109606

$ws.Range("B7").Value = 13086.9856130397
$ws.Range("C7").Value = 12730.5786212774
$ws.Range("E7").Value = 8833.98828694593
$ws.Range("F7").Value = 395.362787842637

$ws.Range("B8").Value = 13086.9856130397
$ws.Range("C8").Value = 12270.4770225147
$ws.Range("E8").Value = 8833.98828694593
$ws.Range("F8").Value = 376.191887894194

$ws.Range("B9").Value = 13086.9856130397
$ws.Range("C9").Value = 12254.8028216669
$ws.Range("E9").Value = 8833.98828694593
$ws.Range("F9").Value = 375.5387961922

$ws.Range("B10").Value = 13086.9856130397
$ws.Range("C10").Value = 11478.3983548796
$ws.Range("E10").Value = 8833.94820637194
$ws.Range("F10").Value = 343.186940052147

$ws.Range("B11").Value = 5313.60020139102
$ws.Range("C11").Value = 7856.97939976793
$ws.Range("E11").Value = 8428.8743126041
$ws.Range("F11").Value = 175.416404682168

$ws.Range("B12").Value = 5204.16189374816
$ws.Range("C12").Value = 7656.6933448276
$ws.Range("E12").Value = 8420.96268051909
$ws.Range("F12").Value = 166.741501056112

$ws.Range("B13").Value = 12789.4588864903
$ws.Range("C13").Value = 11478.959504126
$ws.Range("E13").Value = 8614.55519297089
$ws.Range("F13").Value = 334.068945712369

$ws.Range("B14").Value = 12789.4588864903
$ws.Range("C14").Value = 11741.1213579161
$ws.Range("E14").Value = 8614.55519297089
$ws.Range("F14").Value = 344.992356286957

$ws.Range("B15").Value = 12789.4588864903
$ws.Range("C15").Value = 11908.0066612999
$ws.Range("E15").Value = 8614.55519297089
$ws.Range("F15").Value = 351.945910594617
